$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change B2 from string "1 105 928" to numeric value 1105928
$ws.Range("B2").Value = 1105928

# Update the active cell/selection to B2
$ws.Range("B2").Select()
